$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 7, pushing the existing rows 7-31 down to 9-33.
$ws.Rows("7:8").Insert()

# Fill in row 7 (new data point): Angeleno / Especial, date 2022-03-10
$ws.Range("A7").Value = 2
$ws.Range("B7").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C7").Value = "Coquimbo"
$ws.Range("D7").Value = 44630
$ws.Range("E7").Value = 4
$ws.Range("F7").Value = "Fruta"
$ws.Range("G7").Value = 100103
$ws.Range("H7").Value = "Frutos de hueso (carozo)"
$ws.Range("I7").Value = 100103002
$ws.Range("J7").Value = "Ciruela"
$ws.Range("K7").Value = "Angeleno"
$ws.Range("L7").Value = "Especial"
$ws.Range("M7").Value = 16
$ws.Range("N7").Value = 250000
$ws.Range("O7").Value = 260000
$ws.Range("P7").Value = 255000
$ws.Range("Q7").Value = "$/bins (450 kilos)"
$ws.Range("R7").Value = "Región Metropolitana"
$ws.Range("S7").Value = 567
$ws.Range("T7").Value = 450

# Fill in row 8 (new data point): Angeleno / Primera, date 2022-03-10
$ws.Range("A8").Value = 2
$ws.Range("B8").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C8").Value = "Coquimbo"
$ws.Range("D8").Value = 44630
$ws.Range("E8").Value = 4
$ws.Range("F8").Value = "Fruta"
$ws.Range("G8").Value = 100103
$ws.Range("H8").Value = "Frutos de hueso (carozo)"
$ws.Range("I8").Value = 100103002
$ws.Range("J8").Value = "Ciruela"
$ws.Range("K8").Value = "Angeleno"
$ws.Range("L8").Value = "Primera"
$ws.Range("M8").Value = 20
$ws.Range("N8").Value = 220000
$ws.Range("O8").Value = 230000
$ws.Range("P8").Value = 225000
$ws.Range("Q8").Value = "$/bins (450 kilos)"
$ws.Range("R8").Value = "Región Metropolitana"
$ws.Range("S8").Value = 500
$ws.Range("T8").Value = 450
